$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings so they
# keep their exact formatting (e.g. trailing zeros) as plain text,
# matching the original inline-string cell type.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the crypto price refresh.
$ws.Range("D2").Value = "56.463.39"
$ws.Range("E2").Value = "  -4.19%  "
$ws.Range("D3").Value = "2.400.25"
$ws.Range("E3").Value = "  -4.18%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "501.53"
$ws.Range("E5").Value = "  -6.43%  "
$ws.Range("D6").Value = "128.50"
$ws.Range("E6").Value = "  -4.36%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").Value = "2.397.90"
$ws.Range("E9").Value = "  -4.41%  "
$ws.Range("D10").Value = "0.0956"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("D13").Value = "4.60"
$ws.Range("E13").Value = "  -10.97%  "
$ws.Range("D14").Value = "2.823.80"
$ws.Range("E14").Value = "  -4.27%  "
$ws.Range("D15").Value = "57.044.66"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "21.51"
$ws.Range("E16").Value = "  -3.94%  "
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Value = "2.377.16"
$ws.Range("E18").Value = "  -5.11%  "
$ws.Range("E19").Value = "  -4.97%  "
$ws.Range("D20").Value = "309.86"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("D22").Value = "6.18"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "65.46"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "2.496.12"
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("E27").Value = "  -8.62%  "
$ws.Range("D28").Value = "0.149"
$ws.Range("E28").Value = "  -6.49%  "
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -3.44%  "
$ws.Range("D30").Value = "174.50"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D32").Value = "0.0₃0711"
$ws.Range("E32").Value = "  -6.12%  "
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -7.82%  "
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("E39").Value = "  -5.31%  "
$ws.Range("D40").Value = "35.80"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  -6.24%  "
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("D43").Value = "129.89"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").Value = "3.33"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").Value = "4.77"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("D46").Value = "0.571"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").Value = "253.26"
$ws.Range("E47").Value = "  -7.93%  "
$ws.Range("D48").Value = "0.0897"
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("E49").Value = "  -5.60%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0206"
$ws.Range("E50").Value = "  -5.20%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "16.71"
$ws.Range("E51").Value = "  -4.91%  "
